# Add a new row of data (row 3) to the points sheet, matching the
# existing rows' pattern: name in column A, points in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Gouenji"
$ws.Range("C3").Value = 1
